$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$offices = "0.12% CR/LDUAL+CDM/HBET:6-/Offices`n0.12% CR/LDUAL+CDM/HBET:6-/SOS/Offices`n10.73% CR/LFINF+CDM/H:1/Offices`n1.94% CR/LFINF+CDM/H:2/Offices`n7.80% CR/LFINF+CDM/HBET:3-5/Offices`n1.64% CR/LFINF+CDM/HBET:3-5/SOS/Offices`n5.52% CR+PC/LWAL+CDM/H:1/Offices`n0.59% CR+PC/LWAL+CDM/H:2/Offices`n0.59% CR+PC/LWAL+CDM/HBET:3-5/Offices`n0.35% CR+PC/LWAL+CDM/HBET:6-/Offices`n17.30% MCF/LWAL+CDL/H:1/Offices`n3.40% MCF/LWAL+CDL/H:2/Offices`n6.00% MCF/LWAL+CDL/HBET:3-5/Offices`n41.40% MUR/LWAL+CDN/H:1/Offices`n2.50% MUR/LWAL+CDN/H:2/Offices`n0.0% CR/LFINF+CDL/HBET:3-5/Offices`n0.0% MUR/LWAL+CDN/HBET:3-5/Offices"

$trade = "0.0% CR/LDUAL+CDM/HBET:6-/Trade`n0.0% CR/LDUAL+CDM/HBET:6-/SOS/Trade`n8.73% CR/LFINF+CDM/H:1/Trade`n1.65% CR/LFINF+CDM/H:2/Trade`n0.0% CR/LFINF+CDM/HBET:3-5/Trade`n0.0% CR/LFINF+CDM/HBET:3-5/SOS/Trade`n1.42% CR+PC/LWAL+CDM/H:1/Trade`n0.9% CR+PC/LWAL+CDM/H:2/Trade`n0.0% CR+PC/LWAL+CDM/HBET:3-5/Trade`n0.0% CR+PC/LWAL+CDM/HBET:6-/Trade`n28.8% MCF/LWAL+CDL/H:1/Trade`n3.6% MCF/LWAL+CDL/H:2/Trade`n0.0% MCF/LWAL+CDL/HBET:3-5/Trade`n51.4% MUR/LWAL+CDN/H:1/Trade`n3.5% MUR/LWAL+CDN/H:2/Trade`n0.0% CR/LFINF+CDL/HBET:3-5/Trade`n0.0% MUR/LWAL+CDN/HBET:3-5/Trade"

$hotels = "0.0% CR/LDUAL+CDM/HBET:6-/Hotels`n0.0% CR/LDUAL+CDM/HBET:6-/SOS/Hotels`n6.47% CR/LFINF+CDM/H:1/Hotels`n1.14% CR/LFINF+CDM/H:2/Hotels`n2.72% CR/LFINF+CDM/HBET:3-5/Hotels`n0.31% CR/LFINF+CDM/HBET:3-5/SOS/Hotels`n1.48% CR+PC/LWAL+CDM/H:1/Hotels`n0.34% CR+PC/LWAL+CDM/H:2/Hotels`n0.34% CR+PC/LWAL+CDM/HBET:3-5/Hotels`n0.0% CR+PC/LWAL+CDM/HBET:6-/Hotels`n28.9% MCF/LWAL+CDL/H:1/Hotels`n0.0% MCF/LWAL+CDL/H:2/Hotels`n3.4% MCF/LWAL+CDL/HBET:3-5/Hotels`n51.4% MUR/LWAL+CDN/H:1/Hotels`n1.1% MUR/LWAL+CDN/H:2/Hotels`n0.0% CR/LFINF+CDL/HBET:3-5/Hotels`n2.4% MUR/LWAL+CDN/HBET:3-5/Hotels"

$ws.Range("B2").Value = $offices
$ws.Range("C2").Value = $trade
$ws.Range("D2").Value = $hotels

$ws.Range("C2").Select()
